# Apply the "Updated cryptos list" data refresh (values, URLs, and
# price/volume rows reshuffled per upstream coinranking.com feed).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.756.45'
$ws.Range('E2').Value = '  +5.06%  '
$ws.Range('D3').Value = '2.255.42'
$ws.Range('E3').Value = '  +2.52%  '
$ws.Range('E4').Value = '  +0.12%  '
$ws.Range('D5').Value = "'230.61"
$ws.Range('E5').Value = '  +0.59%  '
$ws.Range('D6').Value = "'0.630"
$ws.Range('E6').Value = '  +2.38%  '
$ws.Range('D7').Value = "'63.14"
$ws.Range('E7').Value = '  +4.72%  '
$ws.Range('E8').Value = '  +0.04%  '
$ws.Range('E9').Value = '  +6.66%  '
$ws.Range('D10').Value = "'0.100"
$ws.Range('E10').Value = '  +12.94%  '
$ws.Range('D11').Value = "'56.24"
$ws.Range('E11').Value = '  -1.23%  '
$ws.Range('B12').Value = 'TRON'
$ws.Range('C12').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D12').Value = "'0.106"
$ws.Range('E12').Value = '  +3.00%  '
$ws.Range('B13').Value = 'Avalanche'
$ws.Range('C13').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D13').Value = "'25.77"
$ws.Range('E13').Value = '  +17.11%  '
$ws.Range('D14').Value = '2.595.28'
$ws.Range('E14').Value = '  +2.51%  '
$ws.Range('D15').Value = "'15.62"
$ws.Range('E15').Value = '  +1.93%  '
$ws.Range('D16').Value = "'5.89"
$ws.Range('E16').Value = '  +6.21%  '
$ws.Range('E17').Value = '  +3.75%  '
$ws.Range('D18').Value = '2.250.13'
$ws.Range('E18').Value = '  +2.08%  '
$ws.Range('D19').Value = '43.705.74'
$ws.Range('E19').Value = '  +4.94%  '
$ws.Range('E20').Value = '  +12.21%  '
$ws.Range('D21').Value = "'73.51"
$ws.Range('E21').Value = '  +2.35%  '
$ws.Range('D22').Value = "'6.02"
$ws.Range('E22').Value = '  +0.05%  '
$ws.Range('D23').Value = "'253.69"
$ws.Range('E23').Value = '  +5.29%  '
$ws.Range('D25').Value = "'2.44"
$ws.Range('E25').Value = '  +4.34%  '
$ws.Range('E26').Value = '  -1.83%  '
$ws.Range('D27').Value = "'9.90"
$ws.Range('E27').Value = '  +3.40%  '
$ws.Range('D28').Value = "'171.32"
$ws.Range('E28').Value = '  +1.91%  '
$ws.Range('B29').Value = 'Kaspa'
$ws.Range('C29').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D29').Value = "'0.137"
$ws.Range('E29').Value = '  -1.36%  '
$ws.Range('B30').Value = 'EthereumClassic'
$ws.Range('C30').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D30').Value = "'20.65"
$ws.Range('E30').Value = '  +5.01%  '
$ws.Range('D31').Value = "'2.83"
$ws.Range('E31').Value = '  +9.29%  '
$ws.Range('E32').Value = '  -3.87%  '
$ws.Range('E33').Value = '  +2.52%  '
$ws.Range('D34').Value = "'0.0676"
$ws.Range('E34').Value = '  +5.17%  '
$ws.Range('D35').Value = "'4.68"
$ws.Range('E35').Value = '  +2.44%  '
$ws.Range('D36').Value = "'4.91"
$ws.Range('E36').Value = '  -0.54%  '
$ws.Range('D37').Value = "'3.87"
$ws.Range('E37').Value = '  +9.15%  '
$ws.Range('D38').Value = "'6.69"
$ws.Range('E38').Value = '  +6.77%  '
$ws.Range('D39').Value = "'2.31"
$ws.Range('E39').Value = '  -0.45%  '
$ws.Range('D40').Value = "'0.0253"
$ws.Range('E40').Value = '  +5.52%  '
$ws.Range('E41').Value = '  +0.30%  '
$ws.Range('D42').Value = "'8.29"
$ws.Range('E42').Value = '  -4.23%  '
$ws.Range('D43').Value = "'17.37"
$ws.Range('E43').Value = '  +8.51%  '
$ws.Range('D44').Value = "'0.0957"
$ws.Range('E44').Value = '  +0.61%  '
$ws.Range('B45').Value = 'TrustWalletToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D45').Value = "'1.19"
$ws.Range('E45').Value = '  -0.46%  '
$ws.Range('B46').Value = 'Aave'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D46').Value = "'96.75"
$ws.Range('E46').Value = '  +0.34%  '
$ws.Range('B47').Value = 'Maker'
$ws.Range('C47').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D47').Value = '1.453.13'
$ws.Range('E47').Value = '  +0.06%  '
$ws.Range('B48').Value = 'FTXToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D48').Value = "'4.30"
$ws.Range('E48').Value = '  -1.19%  '
$ws.Range('B49').Value = 'NEARProtocol'
$ws.Range('C49').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D49').Value = "'2.31"
$ws.Range('E49').Value = '  +5.12%  '
$ws.Range('B50').Value = 'TerraClassic'
$ws.Range('C50').Value = 'https://coinranking.com/coin/AaQUAs2Mc+terraclassic-lunc'
$ws.Range('D50').Value = "'0.000205"
$ws.Range('E50').Value = '  -14.86%  '
$ws.Range('E51').Value = '  +0.55%  '

Write-Output "Updated 108 cells"
